$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-14 switch from style index 3 (italic-ish "Normal 2") back to style index 2
# ("Normal 3") - same style used by the surrounding rows (11, 15, 16).
$ws.Range("A12:N14").Style = $ws.Range("A11:N11").Style

# New row 17: "Abandonded Chapel"
$ws.Range("A17").Value = "Abandonded Chapel"
$ws.Range("E17").Value = "An old decrepid chapel in the middle of no where. Half burned, half rotted, what remains is a story of the past."
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 208
$ws.Range("K17").Value = 416
$ws.Range("M17").Value = "Yes"

# Row 17 keeps the style that rows 12-14 originally had (style index 3).
$ws.Range("A17:N17").Style = $ws.Range("A12:N12").Style

$ws.Range("A17").Select()
